# Edit the "flow" slide (slide 3) of the presentation:
#  - "Hash all of the code files" -> "Hash all of the code blocks"
#  - "Compare hash code for each file with every other one"
#       -> "Compare hash code for each block with every other one"
# Both edits are applied run-by-run (matching how PowerPoint splits runs
# when text is retyped in place) so the resulting OOXML has the same
# paragraph/run structure as the authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Shape 5 ("TextBox 7"): Hash all of the code files -> ... blocks ---
$shp1 = $s.Shapes.Item(5)
$tr1 = $shp1.TextFrame.TextRange
$tr1.Text = "Hash all of the code "
$run1b = $tr1.InsertAfter("blocks")

# --- Shape 6 ("TextBox 3"): Compare hash code for each file with every other one ---
$shp2 = $s.Shapes.Item(6)
$tr2 = $shp2.TextFrame.TextRange
$tr2.Text = "Compare hash code for "
$run2b = $tr2.InsertAfter("each ")
$run2c = $run2b.InsertAfter("block with ")
$run2d = $run2c.InsertAfter("every other one")
